$wb = $excel.ActiveWorkbook

# Sheet: Neodymium
$ws = $wb.Worksheets.Item("Neodymium")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = 0.000002195405251500087
$ws.Range("C3").Value = 0.0001062411525673284
$ws.Range("C4").Value = 0.00009608716352691787
$ws.Range("C5").Value = 0.000000002138791829054013

# Sheet: Dysprosium
$ws = $wb.Worksheets.Item("Dysprosium")
$ws.Range("C1").Value = 2030

# Sheet: Copper
$ws = $wb.Worksheets.Item("Copper")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = 0.003816340722347758
$ws.Range("C3").Value = 0.0137679456486295
$ws.Range("C4").Value = 0.003685389348936462
$ws.Range("C5").Value = 0.008080293662215502

# Sheet: Raw silicon
$ws = $wb.Worksheets.Item("Raw silicon")
$ws.Range("C1").Value = 2030
$ws.Range("C2").Value = 0.00005750015024097242
$ws.Range("C3").Value = 0.0001921210602835477
$ws.Range("C4").Value = 0.000053924808017845
$ws.Range("C5").Value = 0.00006847896595910317
